# Fipe_temp.xlsx edit: add Aston Martin "Rapide/Vanquish/Vantage" model rows and
# refresh the dependent Ano/CodigoFipe/PrecoMedio/Mes columns so the shared
# string table stays consistent (commit: "Adaptando para caminhoes e motos").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the "PrecoMedio" column as text (values look numeric, e.g. " 27313.00",
# and must keep their leading space / trailing zeros instead of being coerced
# into real numbers by Excel's type inference).
$ws.Range("E2:E183").NumberFormat = "@"

# Rows 2-164 already existed; only columns C (AnoSelecionado), D (CodigoFipe),
# E (PrecoMedio) and F (Mes Referencia) shift to new shared-string values.
$existingChanges = @"
2|3|1999 Gasolina
2|4|006009-7
2|5| 27313.00
2|6|junho de 2025
3|3|1998 Gasolina
3|4|006009-7
3|5| 21558.00
3|6|junho de 2025
4|3|1997 Gasolina
4|4|006009-7
4|5| 14581.00
4|6|junho de 2025
5|3|1996 Gasolina
5|4|006009-7
5|5| 13265.00
5|6|junho de 2025
6|3|1998 Gasolina
6|4|006001-1
6|5| 21581.00
6|6|junho de 2025
7|3|1997 Gasolina
7|4|006001-1
7|5| 17680.00
7|6|junho de 2025
8|3|1996 Gasolina
8|4|006001-1
8|5| 17248.00
8|6|junho de 2025
9|3|1995 Gasolina
9|4|006001-1
9|5| 13550.00
9|6|junho de 2025
10|3|1999 Gasolina
10|4|006002-0
10|5| 29872.00
10|6|junho de 2025
11|3|1998 Gasolina
11|4|006002-0
11|5| 29143.00
11|6|junho de 2025
12|3|1997 Gasolina
12|4|006002-0
12|5| 23404.00
12|6|junho de 2025
13|3|1996 Gasolina
13|4|006002-0
13|5| 22833.00
13|6|junho de 2025
14|3|1998 Gasolina
14|4|006008-9
14|5| 22449.00
14|6|junho de 2025
15|3|1997 Gasolina
15|4|006008-9
15|5| 18647.00
15|6|junho de 2025
16|3|1996 Gasolina
16|4|006008-9
16|5| 18192.00
16|6|junho de 2025
17|3|2005 Gasolina
17|4|006017-8
17|5| 92142.00
17|6|junho de 2025
18|3|2004 Gasolina
18|4|006017-8
18|5| 89894.00
18|6|junho de 2025
19|3|2003 Gasolina
19|4|006017-8
19|5| 87701.00
19|6|junho de 2025
20|3|1997 Gasolina
20|4|006003-8
20|5| 19511.00
20|6|junho de 2025
21|3|1996 Gasolina
21|4|006003-8
21|5| 18466.00
21|6|junho de 2025
22|3|1995 Gasolina
22|4|006003-8
22|5| 15931.00
22|6|junho de 2025
23|3|1997 Gasolina
23|4|006004-6
23|5| 24925.00
23|6|junho de 2025
24|3|1996 Gasolina
24|4|006004-6
24|5| 24317.00
24|6|junho de 2025
25|3|1995 Gasolina
25|4|006004-6
25|5| 19754.00
25|6|junho de 2025
26|3|2003 Gasolina
26|4|006015-1
26|5| 104853.00
26|6|junho de 2025
27|3|1999 Gasolina
27|4|006015-1
27|5| 32424.00
27|6|junho de 2025
28|3|2002 Gasolina
28|4|006014-3
28|5| 49978.00
28|6|junho de 2025
29|3|2001 Gasolina
29|4|006014-3
29|5| 45101.00
29|6|junho de 2025
30|3|2004 Gasolina
30|4|006016-0
30|5| 74355.00
30|6|junho de 2025
31|3|2003 Gasolina
31|4|006016-0
31|5| 71998.00
31|6|junho de 2025
32|3|2002 Gasolina
32|4|006010-0
32|5| 75561.00
32|6|junho de 2025
33|3|2001 Gasolina
33|4|006010-0
33|5| 70007.00
33|6|junho de 2025
34|3|2000 Gasolina
34|4|006010-0
34|5| 58514.00
34|6|junho de 2025
35|3|1999 Gasolina
35|4|006010-0
35|5| 57086.00
35|6|junho de 2025
36|3|1996 Gasolina
36|4|006005-4
36|5| 29151.00
36|6|junho de 2025
37|3|1995 Gasolina
37|4|006005-4
37|5| 28440.00
37|6|junho de 2025
38|3|1994 Gasolina
38|4|006005-4
38|5| 27102.00
38|6|junho de 2025
39|3|1993 Gasolina
39|4|006005-4
39|5| 19575.00
39|6|junho de 2025
40|3|1992 Gasolina
40|4|006005-4
40|5| 13455.00
40|6|junho de 2025
41|3|1991 Gasolina
41|4|006005-4
41|5| 11701.00
41|6|junho de 2025
42|3|1996 Gasolina
42|4|006006-2
42|5| 28869.00
42|6|junho de 2025
43|3|1995 Gasolina
43|4|006006-2
43|5| 28164.00
43|6|junho de 2025
44|3|2003 Gasolina
44|4|006011-9
44|5| 89797.00
44|6|junho de 2025
45|3|2002 Gasolina
45|4|006011-9
45|5| 86004.00
45|6|junho de 2025
46|3|2001 Gasolina
46|4|006011-9
46|5| 83906.00
46|6|junho de 2025
47|3|2000 Gasolina
47|4|006011-9
47|5| 81859.00
47|6|junho de 2025
48|3|1999 Gasolina
48|4|006011-9
48|5| 71424.00
48|6|junho de 2025
49|3|1986 Gasolina
49|4|006013-5
49|5| 8070.00
49|6|junho de 2025
50|3|1985 Gasolina
50|4|006013-5
50|5| 7072.00
50|6|junho de 2025
51|3|1997 Gasolina
51|4|006007-0
51|5| 133185.00
51|6|junho de 2025
52|3|1996 Gasolina
52|4|006007-0
52|5| 129936.00
52|6|junho de 2025
53|3|1995 Gasolina
53|4|006007-0
53|5| 92099.00
53|6|junho de 2025
54|3|1993 Gasolina
54|4|006007-0
54|5| 41220.00
54|6|junho de 2025
55|3|1992 Gasolina
55|4|006007-0
55|5| 32289.00
55|6|junho de 2025
56|3|1991 Gasolina
56|4|006007-0
56|5| 27713.00
56|6|junho de 2025
57|3|2000 Diesel
57|4|037001-0
57|5| 426859.00
57|6|junho de 2025
58|3|1999 Diesel
58|4|037001-0
58|5| 384365.00
58|6|junho de 2025
59|3|1998 Diesel
59|4|037001-0
59|5| 332825.00
59|6|junho de 2025
60|3|2000 Diesel
60|4|037002-9
60|5| 396648.00
60|6|junho de 2025
61|3|1999 Diesel
61|4|037002-9
61|5| 351603.00
61|6|junho de 2025
62|3|1998 Diesel
62|4|037002-9
62|5| 294624.00
62|6|junho de 2025
63|3|2000 Diesel
63|4|037003-7
63|5| 453297.00
63|6|junho de 2025
64|3|1999 Diesel
64|4|037003-7
64|5| 401677.00
64|6|junho de 2025
65|3|1998 Diesel
65|4|037003-7
65|5| 348246.00
65|6|junho de 2025
66|3|1998 Diesel
66|4|007018-1
66|5| 25397.00
66|6|junho de 2025
67|3|1997 Diesel
67|4|007018-1
67|5| 21995.00
67|6|junho de 2025
68|3|1996 Diesel
68|4|007018-1
68|5| 18798.00
68|6|junho de 2025
69|3|1995 Diesel
69|4|007018-1
69|5| 18158.00
69|6|junho de 2025
70|3|1994 Diesel
70|4|007018-1
70|5| 16980.00
70|6|junho de 2025
71|3|1998 Diesel
71|4|007019-0
71|5| 27427.00
71|6|junho de 2025
72|3|1997 Diesel
72|4|007019-0
72|5| 25605.00
72|6|junho de 2025
73|3|1996 Diesel
73|4|007019-0
73|5| 21492.00
73|6|junho de 2025
74|3|1995 Diesel
74|4|007019-0
74|5| 19903.00
74|6|junho de 2025
75|3|1994 Diesel
75|4|007019-0
75|5| 19231.00
75|6|junho de 2025
76|3|1998 Diesel
76|4|007001-7
76|5| 15514.00
76|6|junho de 2025
77|3|1997 Diesel
77|4|007001-7
77|5| 14842.00
77|6|junho de 2025
78|3|1996 Diesel
78|4|007001-7
78|5| 12659.00
78|6|junho de 2025
79|3|1995 Diesel
79|4|007001-7
79|5| 12058.00
79|6|junho de 2025
80|3|1994 Diesel
80|4|007001-7
80|5| 11358.00
80|6|junho de 2025
81|3|1993 Diesel
81|4|007001-7
81|5| 9903.00
81|6|junho de 2025
82|3|1998 Diesel
82|4|007002-5
82|5| 14608.00
82|6|junho de 2025
83|3|1997 Diesel
83|4|007002-5
83|5| 13421.00
83|6|junho de 2025
84|3|1996 Diesel
84|4|007002-5
84|5| 12439.00
84|6|junho de 2025
85|3|1995 Diesel
85|4|007002-5
85|5| 11373.00
85|6|junho de 2025
86|3|1994 Diesel
86|4|007002-5
86|5| 10192.00
86|6|junho de 2025
87|3|1993 Diesel
87|4|007002-5
87|5| 9297.00
87|6|junho de 2025
88|3|1997 Diesel
88|4|007017-3
88|5| 12164.00
88|6|junho de 2025
89|3|1996 Diesel
89|4|007017-3
89|5| 11236.00
89|6|junho de 2025
90|3|1995 Diesel
90|4|007017-3
90|5| 9800.00
90|6|junho de 2025
91|3|1994 Diesel
91|4|007017-3
91|5| 9324.00
91|6|junho de 2025
92|3|1995 Diesel
92|4|007020-3
92|5| 25330.00
92|6|junho de 2025
93|3|1994 Diesel
93|4|007020-3
93|5| 24712.00
93|6|junho de 2025
94|3|1999 Diesel
94|4|007016-5
94|5| 16727.00
94|6|junho de 2025
95|3|1998 Diesel
95|4|007016-5
95|5| 13972.00
95|6|junho de 2025
96|3|1997 Diesel
96|4|007016-5
96|5| 13585.00
96|6|junho de 2025
97|3|1999 Diesel
97|4|007014-9
97|5| 16655.00
97|6|junho de 2025
98|3|1998 Diesel
98|4|007014-9
98|5| 14777.00
98|6|junho de 2025
99|3|1997 Diesel
99|4|007014-9
99|5| 12873.00
99|6|junho de 2025
100|3|1999 Diesel
100|4|007015-7
100|5| 19837.00
100|6|junho de 2025
101|3|1998 Diesel
101|4|007016-5
101|5| 13972.00
101|6|junho de 2025
102|3|1997 Diesel
102|4|007015-7
102|5| 14311.00
102|6|junho de 2025
103|3|1996 Diesel
103|4|007015-7
103|5| 12443.00
103|6|junho de 2025
104|3|1998 Gasolina
104|4|007003-3
104|5| 6442.00
104|6|junho de 2025
105|3|1997 Gasolina
105|4|007003-3
105|5| 6095.00
105|6|junho de 2025
106|3|1996 Gasolina
106|4|007003-3
106|5| 5884.00
106|6|junho de 2025
107|3|1995 Gasolina
107|4|007003-3
107|5| 5680.00
107|6|junho de 2025
108|3|1999 Gasolina
108|4|007009-2
108|5| 5955.00
108|6|junho de 2025
109|3|1998 Gasolina
109|4|007009-2
109|5| 5298.00
109|6|junho de 2025
110|3|1997 Gasolina
110|4|007009-2
110|5| 5168.00
110|6|junho de 2025
111|3|1996 Gasolina
111|4|007009-2
111|5| 4546.00
111|6|junho de 2025
112|3|1995 Gasolina
112|4|007009-2
112|5| 4415.00
112|6|junho de 2025
113|3|1994 Gasolina
113|4|007009-2
113|5| 3607.00
113|6|junho de 2025
114|3|1998 Gasolina
114|4|007004-1
114|5| 5607.00
114|6|junho de 2025
115|3|1997 Gasolina
115|4|007004-1
115|5| 5281.00
115|6|junho de 2025
116|3|1996 Gasolina
116|4|007004-1
116|5| 4945.00
116|6|junho de 2025
117|3|1995 Gasolina
117|4|007004-1
117|5| 4532.00
117|6|junho de 2025
118|3|1994 Gasolina
118|4|007004-1
118|5| 3838.00
118|6|junho de 2025
119|3|1993 Gasolina
119|4|007004-1
119|5| 3291.00
119|6|junho de 2025
120|3|1999 Gasolina
120|4|007012-2
120|5| 5858.00
120|6|junho de 2025
121|3|1998 Gasolina
121|4|007012-2
121|5| 5699.00
121|6|junho de 2025
122|3|1997 Gasolina
122|4|007012-2
122|5| 5481.00
122|6|junho de 2025
123|3|1999 Gasolina
123|4|007011-4
123|5| 4160.00
123|6|junho de 2025
124|3|1998 Gasolina
124|4|007011-4
124|5| 4058.00
124|6|junho de 2025
125|3|1997 Gasolina
125|4|007011-4
125|5| 3746.00
125|6|junho de 2025
126|3|1998 Gasolina
126|4|007005-0
126|5| 3569.00
126|6|junho de 2025
127|3|1997 Gasolina
127|4|007005-0
127|5| 3442.00
127|6|junho de 2025
128|3|1996 Gasolina
128|4|007005-0
128|5| 3284.00
128|6|junho de 2025
129|3|1995 Gasolina
129|4|007005-0
129|5| 2840.00
129|6|junho de 2025
130|3|1994 Gasolina
130|4|007005-0
130|5| 2401.00
130|6|junho de 2025
131|3|1999 Gasolina
131|4|007008-4
131|5| 5876.00
131|6|junho de 2025
132|3|1998 Gasolina
132|4|007008-4
132|5| 5732.00
132|6|junho de 2025
133|3|1997 Gasolina
133|4|007008-4
133|5| 4374.00
133|6|junho de 2025
134|3|1998 Gasolina
134|4|007006-8
134|5| 7989.00
134|6|junho de 2025
135|3|1997 Gasolina
135|4|007006-8
135|5| 7146.00
135|6|junho de 2025
136|3|1996 Gasolina
136|4|007006-8
136|5| 6889.00
136|6|junho de 2025
137|3|1995 Gasolina
137|4|007006-8
137|5| 6037.00
137|6|junho de 2025
138|3|1994 Gasolina
138|4|007006-8
138|5| 4876.00
138|6|junho de 2025
139|3|1993 Gasolina
139|4|007006-8
139|5| 4281.00
139|6|junho de 2025
140|3|1999 Gasolina
140|4|007013-0
140|5| 6241.00
140|6|junho de 2025
141|3|1998 Gasolina
141|4|007013-0
141|5| 6088.00
141|6|junho de 2025
142|3|1997 Gasolina
142|4|007013-0
142|5| 5939.00
142|6|junho de 2025
143|3|1998 Gasolina
143|4|007007-6
143|5| 5967.00
143|6|junho de 2025
144|3|1997 Gasolina
144|4|007007-6
144|5| 5821.00
144|6|junho de 2025
145|3|1996 Gasolina
145|4|007007-6
145|5| 5569.00
145|6|junho de 2025
146|3|1995 Gasolina
146|4|007007-6
146|5| 4635.00
146|6|junho de 2025
147|3|1994 Gasolina
147|4|007007-6
147|5| 4206.00
147|6|junho de 2025
148|3|1993 Gasolina
148|4|007007-6
148|5| 4005.00
148|6|junho de 2025
149|3|Zero KM a Gasolina
149|4|085018-7
149|5| 3600000.00
149|6|junho de 2025
150|3|2025 Gasolina
150|4|085018-7
150|5| 3228897.00
150|6|junho de 2025
151|3|2024 Gasolina
151|4|085018-7
151|5| 3070680.00
151|6|junho de 2025
152|3|Zero KM a Gasolina
152|4|085019-5
152|5| 3900000.00
152|6|junho de 2025
153|3|2025 Gasolina
153|4|085019-5
153|5| 3625334.00
153|6|junho de 2025
154|3|2016 Gasolina
154|4|085011-0
154|5| 1214247.00
154|6|junho de 2025
155|3|2014 Gasolina
155|4|085011-0
155|5| 946995.00
155|6|junho de 2025
156|3|2011 Gasolina
156|4|085006-3
156|5| 617791.00
156|6|junho de 2025
157|3|2023 Gasolina
157|4|085014-4
157|5| 4493874.00
157|6|junho de 2025
158|3|2022 Gasolina
158|4|085014-4
158|5| 4142252.00
158|6|junho de 2025
159|3|2023 Gasolina
159|4|085015-2
159|5| 2346605.00
159|6|junho de 2025
160|3|2022 Gasolina
160|4|085015-2
160|5| 2055673.00
160|6|junho de 2025
161|3|Zero KM a Gasolina
161|4|085016-0
161|5| 3800000.00
161|6|junho de 2025
162|3|2025 Gasolina
162|4|085016-0
162|5| 3226068.00
162|6|junho de 2025
163|3|2024 Gasolina
163|4|085016-0
163|5| 2729902.00
163|6|junho de 2025
164|3|2023 Gasolina
164|4|085016-0
164|5| 2502679.00
164|6|junho de 2025
"@

foreach ($line in ($existingChanges -split "`n")) {
    $line = $line.Trim("`r", "`n")
    if ($line -eq "") { continue }
    $parts = $line -split '\|', 3
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $val = $parts[2]
    $ws.Cells.Item($r, $c).Value = $val
}

# Rows 165-183 are brand new ASTON MARTIN entries (Rapide, Rapide S, Vanquish,
# Vantage variants) appended after the existing DBX707 rows.
$newRows = @"
165|1|ASTON MARTIN
165|2|Rapide 6.0 V12 477cv
165|3|2012 Gasolina
165|4|085007-1
165|5| 694217.00
165|6|junho de 2025
166|1|ASTON MARTIN
166|2|Rapide 6.0 V12 477cv
166|3|2011 Gasolina
166|4|085007-1
166|5| 607011.00
166|6|junho de 2025
167|1|ASTON MARTIN
167|2|Rapide S 6.0 V12 550cv
167|3|2016 Gasolina
167|4|085009-8
167|5| 1296696.00
167|6|junho de 2025
168|1|ASTON MARTIN
168|2|Rapide S 6.0 V12 550cv
168|3|2014 Gasolina
168|4|085009-8
168|5| 1032188.00
168|6|junho de 2025
169|1|ASTON MARTIN
169|2|Vanquish V12 6.0 565cv
169|3|2016 Gasolina
169|4|085010-1
169|5| 1784133.00
169|6|junho de 2025
170|1|ASTON MARTIN
170|2|Vanquish V12 6.0 565cv
170|3|2014 Gasolina
170|4|085010-1
170|5| 1616174.00
170|6|junho de 2025
171|1|ASTON MARTIN
171|2|Vantage 6.0 V12 510cv
171|3|2012 Gasolina
171|4|085004-7
171|5| 476063.00
171|6|junho de 2025
172|1|ASTON MARTIN
172|2|Vantage Coupe 4.7 V8 425cv
172|3|2016 Gasolina
172|4|085002-0
172|5| 638283.00
172|6|junho de 2025
173|1|ASTON MARTIN
173|2|Vantage Coupe 4.7 V8 425cv
173|3|2014 Gasolina
173|4|085002-0
173|5| 591028.00
173|6|junho de 2025
174|1|ASTON MARTIN
174|2|Vantage Coupe 4.7 V8 425cv
174|3|2012 Gasolina
174|4|085002-0
174|5| 528207.00
174|6|junho de 2025
175|1|ASTON MARTIN
175|2|Vantage Coupe 4.7 V8 425cv
175|3|2011 Gasolina
175|4|085002-0
175|5| 515323.00
175|6|junho de 2025
176|1|ASTON MARTIN
176|2|Vantage Cupê  4.0 V8 510cv
176|3|Zero KM a Gasolina
176|4|085012-8
176|5| 2800750.00
176|6|junho de 2025
177|1|ASTON MARTIN
177|2|Vantage Cupê  4.0 V8 510cv
177|3|2025 Gasolina
177|4|085012-8
177|5| 2200329.00
177|6|junho de 2025
178|1|ASTON MARTIN
178|2|Vantage Cupê  4.0 V8 510cv
178|3|2023 Gasolina
178|4|085012-8
178|5| 1739268.00
178|6|junho de 2025
179|1|ASTON MARTIN
179|2|Vantage Cupê  4.0 V8 510cv
179|3|2022 Gasolina
179|4|085012-8
179|5| 1647750.00
179|6|junho de 2025
180|1|ASTON MARTIN
180|2|Vantage Cupê F1 Edition 4.0 V8 535cv
180|3|2023 Gasolina
180|4|085013-6
180|5| 1770690.00
180|6|junho de 2025
181|1|ASTON MARTIN
181|2|Vantage Cupê F1 Edition 4.0 V8 535cv
181|3|2022 Gasolina
181|4|085013-6
181|5| 1706250.00
181|6|junho de 2025
182|1|ASTON MARTIN
182|2|Vantage Roadster 4.7 V8 420cv
182|3|2011 Gasolina
182|4|085003-9
182|5| 588516.00
182|6|junho de 2025
183|1|ASTON MARTIN
183|2|Vantage Roadster 4.7 V8 420cv
183|3|2010 Gasolina
183|4|085003-9
183|5| 521625.00
183|6|junho de 2025
"@

foreach ($line in ($newRows -split "`n")) {
    $line = $line.Trim("`r", "`n")
    if ($line -eq "") { continue }
    $parts = $line -split '\|', 3
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $val = $parts[2]
    $ws.Cells.Item($r, $c).Value = $val
}

Write-Output "Edit complete"
